# Auto-generated script to apply Chocobo_Profits.xlsx market-data refresh
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ,
# LevePriceNQ / LevePriceHQ and the derived LeveProfitNQ / LeveProfitHQ columns
# (H:N) for the rows whose cached market-board data changed.

$wb = $excel.ActiveWorkbook

# --- Worksheet "ALC" ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 520.8182
$ws.Range("I39").Value = 158.16667
$ws.Range("J39").Value = 956
$ws.Range("K39").Value = 474.50001
$ws.Range("L39").Value = 2868
$ws.Range("M39").Value = -178.50001
$ws.Range("N39").Value = -3460
$ws.Range("H80").Value = 981.3125
$ws.Range("I80").Value = 362.625
$ws.Range("J80").Value = 1600
$ws.Range("K80").Value = 1087.875
$ws.Range("L80").Value = 4800
$ws.Range("M80").Value = -89.875
$ws.Range("N80").Value = -6796
$ws.Range("H83").Value = 981.3125
$ws.Range("I83").Value = 362.625
$ws.Range("J83").Value = 1600
$ws.Range("K83").Value = 3263.625
$ws.Range("L83").Value = 14400
$ws.Range("M83").Value = 1728.375
$ws.Range("N83").Value = -24384
$ws.Range("H116").Value = 405822.44
$ws.Range("I116").Value = 1430414.2
$ws.Range("J116").Value = 7370.0557
$ws.Range("K116").Value = 1430414.2
$ws.Range("L116").Value = 7370.0557
$ws.Range("M116").Value = -1426972.2
$ws.Range("N116").Value = -14254.0557
$ws.Range("H135").Value = 540.7
$ws.Range("I135").Value = 490.21054
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 4411.894859999999
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -1876.894859999999
$ws.Range("N135").Value = -18570
$ws.Range("H137").Value = 1907238.6
$ws.Range("I137").Value = 2647114.5
$ws.Range("J137").Value = 4700.4287
$ws.Range("K137").Value = 7941343.5
$ws.Range("L137").Value = 14101.2861
$ws.Range("M137").Value = -7938793.5
$ws.Range("N137").Value = -19201.2861
$ws.Range("H141").Value = 12367.1
$ws.Range("I141").Value = 16124.429
$ws.Range("J141").Value = 3600
$ws.Range("K141").Value = 48373.287
$ws.Range("L141").Value = 10800
$ws.Range("M141").Value = -43193.287
$ws.Range("N141").Value = -21160

# --- Worksheet "ARM" ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1573.27
$ws.Range("I32").Value = 1291.7922
$ws.Range("J32").Value = 2515.6086
$ws.Range("K32").Value = 1291.7922
$ws.Range("L32").Value = 2515.6086
$ws.Range("M32").Value = -1004.7922
$ws.Range("N32").Value = -3089.6086
$ws.Range("H88").Value = 6670075
$ws.Range("I88").Value = 33334232
$ws.Range("J88").Value = 4035.875
$ws.Range("K88").Value = 33334232
$ws.Range("L88").Value = 4035.875
$ws.Range("M88").Value = -33333826
$ws.Range("N88").Value = -4847.875
$ws.Range("H91").Value = 6670075
$ws.Range("I91").Value = 33334232
$ws.Range("J91").Value = 4035.875
$ws.Range("K91").Value = 33334232
$ws.Range("L91").Value = 4035.875
$ws.Range("M91").Value = -33332828
$ws.Range("N91").Value = -6843.875
$ws.Range("H109").Value = 31594.875
$ws.Range("J109").Value = 31594.875
$ws.Range("L109").Value = 31594.875
$ws.Range("N109").Value = -34368.875
$ws.Range("H122").Value = 3502.8367
$ws.Range("I122").Value = 2949.9211
$ws.Range("J122").Value = 5412.909
$ws.Range("K122").Value = 8849.763300000001
$ws.Range("L122").Value = 16238.727
$ws.Range("M122").Value = -6399.763300000001
$ws.Range("N122").Value = -21138.727
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null
$ws.Range("H137").Value = 45780
$ws.Range("J137").Value = 45780
$ws.Range("L137").Value = 45780
$ws.Range("N137").Value = -55980

# --- Worksheet "BSM" ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H114").Value = 39999
$ws.Range("J114").Value = 39999
$ws.Range("L114").Value = 39999
$ws.Range("N114").Value = -48677
$ws.Range("H137").Value = 48920
$ws.Range("J137").Value = 48920
$ws.Range("L137").Value = 48920
$ws.Range("N137").Value = -59120

# --- Worksheet "CRP" ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11629.211
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 11629.211
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 11629.211
$ws.Range("M31").Value = $null
$ws.Range("N31").Value = -12219.211
$ws.Range("H34").Value = 11629.211
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 11629.211
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 11629.211
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = -12033.211
$ws.Range("H98").Value = 47199
$ws.Range("J98").Value = 47199
$ws.Range("L98").Value = 47199
$ws.Range("N98").Value = -51691
$ws.Range("H99").Value = 14290193
$ws.Range("I99").Value = 40001660
$ws.Range("J99").Value = 6044.4443
$ws.Range("K99").Value = 40001660
$ws.Range("L99").Value = 6044.4443
$ws.Range("M99").Value = -40000162
$ws.Range("N99").Value = -9040.444299999999
$ws.Range("H122").Value = 3828.4285
$ws.Range("I122").Value = 1799
$ws.Range("J122").Value = 4166.6665
$ws.Range("K122").Value = 5397
$ws.Range("L122").Value = 12499.9995
$ws.Range("M122").Value = -2947
$ws.Range("N122").Value = -17399.9995
$ws.Range("H126").Value = 14290193
$ws.Range("I126").Value = 40001660
$ws.Range("J126").Value = 6044.4443
$ws.Range("K126").Value = 120004980
$ws.Range("L126").Value = 18133.3329
$ws.Range("M126").Value = -120002510
$ws.Range("N126").Value = -23073.3329
$ws.Range("H137").Value = 30480
$ws.Range("J137").Value = 30480
$ws.Range("L137").Value = 30480
$ws.Range("N137").Value = -40680

# --- Worksheet "CUL" ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 290.2
$ws.Range("I46").Value = 290.2
$ws.Range("K46").Value = 870.5999999999999
$ws.Range("M46").Value = -779.5999999999999
$ws.Range("H113").Value = 4808260.5
$ws.Range("J113").Value = 11364175
$ws.Range("L113").Value = 34092525
$ws.Range("N113").Value = -34096865
$ws.Range("H131").Value = 787.5700000000001
$ws.Range("I131").Value = 305.44446
$ws.Range("J131").Value = 835.25275
$ws.Range("K131").Value = 916.33338
$ws.Range("L131").Value = 2505.75825
$ws.Range("M131").Value = 4123.66662
$ws.Range("N131").Value = -12585.75825
$ws.Range("H132").Value = 1848.4166
$ws.Range("I132").Value = 766.3333
$ws.Range("K132").Value = 6896.9997
$ws.Range("M132").Value = -4366.9997

# --- Worksheet "GSM" ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 27817.455
$ws.Range("J46").Value = 28074.2
$ws.Range("L46").Value = 28074.2
$ws.Range("N46").Value = -28386.2
$ws.Range("H126").Value = 3140.8
$ws.Range("I126").Value = 2722.078
$ws.Range("J126").Value = 4542.609
$ws.Range("K126").Value = 8166.234
$ws.Range("L126").Value = 13627.827
$ws.Range("M126").Value = -5696.234
$ws.Range("N126").Value = -18567.827
$ws.Range("H137").Value = 43750
$ws.Range("J137").Value = 43750
$ws.Range("L137").Value = 43750
$ws.Range("N137").Value = -53950

# --- Worksheet "LTW" ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1174.2258
$ws.Range("J68").Value = 3156.6667
$ws.Range("L68").Value = 3156.6667
$ws.Range("N68").Value = -4654.6667
$ws.Range("H71").Value = 1174.2258
$ws.Range("J71").Value = 3156.6667
$ws.Range("L71").Value = 15783.3335
$ws.Range("N71").Value = -23271.3335
$ws.Range("H132").Value = 8732.1875
$ws.Range("I132").Value = 4890.222
$ws.Range("J132").Value = 13671.857
$ws.Range("K132").Value = 14670.666
$ws.Range("L132").Value = 41015.571
$ws.Range("M132").Value = -12140.666
$ws.Range("N132").Value = -46075.571

# --- Worksheet "WVR" ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 18250
$ws.Range("J57").Value = 18250
$ws.Range("L57").Value = 18250
$ws.Range("N57").Value = -19758
$ws.Range("H80").Value = 39833.332
$ws.Range("J80").Value = 39833.332
$ws.Range("L80").Value = 39833.332
$ws.Range("N80").Value = -41829.332
$ws.Range("H83").Value = 39833.332
$ws.Range("J83").Value = 39833.332
$ws.Range("L83").Value = 119499.996
$ws.Range("N83").Value = -129483.996
$ws.Range("H123").Value = 34938.75
$ws.Range("J123").Value = 34938.75
$ws.Range("L123").Value = 34938.75
$ws.Range("N123").Value = -44738.75
